# Generate Report for Handback
# Re-populate the localization-status workbook so the c9f7cf44-... file
# (previously "Ready for handoff") now shows as handed-back, and reorder
# the per-file rows so the most recently processed file (c9f7cf44) is
# listed first on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-03-24 19:28:46"

$ov.Range("A3").Value = "ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"
$ov.Range("D3").Value = "2016-03-24 19:26:54"

$ov.Range("A4").Value = "ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md"
$ov.Range("B4").Value = "Handed back: in sync with en-US"
$ov.Range("C4").Value = "Handed back: in sync with en-US"
$ov.Range("D4").Value = "2016-03-24 19:26:54"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md", "", "", "ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md", "", "", "ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2: c9f7cf44 file, now handed back
$zh.Range("A2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-24 19:28:42"
$zh.Range("F2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md"
$zh.Range("G2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-24 19:29:08"
$zh.Range("J2").Value = "Include"

# Row 3: ffff86e48609 file
$zh.Range("A3").Value = "ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-24 19:26:49"
$zh.Range("F3").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.md"
$zh.Range("G3").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-24 19:27:14"
$zh.Range("J3").Value = "Include"

# Row 4: ffffff7695d231 file
$zh.Range("A4").Value = "ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("D4").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-24 19:26:49"
$zh.Range("F4").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.md"
$zh.Range("G4").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf"
$zh.Range("H4").Value = "2016-03-24 19:27:14"
$zh.Range("J4").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ac11cc8578733d63aecbe27859039e3510f1f36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.zh-cn.xlf", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ac11cc8578733d63aecbe27859039e3510f1f36/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.zh-cn.xlf", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md", "", "", "ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d80b767201cb984e276df16d9397214d1f05ff5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/258372e0c314f54d7f73eac1563554cfc2184482/e2e/ee944fac-2a76-4622-b02c-316995cd2ba5.md", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc1b70a45685793f01477a5274d0297a5b93d735/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf")

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md", "", "", "ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d80b767201cb984e276df16d9397214d1f05ff5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/258372e0c314f54d7f73eac1563554cfc2184482/e2e/ee944fac-2a76-4622-b02c-316995cd2ba5.md", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.md")
$zh.Hyperlinks.Add($zh.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc1b70a45685793f01477a5274d0297a5b93d735/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2: c9f7cf44 file, now handed back
$de.Range("A2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.de-de.xlf"
$de.Range("E2").Value = "2016-03-24 19:28:46"
$de.Range("F2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md"
$de.Range("G2").Value = "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.de-de.xlf"
$de.Range("H2").Value = "2016-03-24 19:29:15"
$de.Range("J2").Value = "Include"

# Row 3: ffff86e48609 file
$de.Range("A3").Value = "ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf"
$de.Range("E3").Value = "2016-03-24 19:26:54"
$de.Range("F3").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.md"
$de.Range("G3").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf"
$de.Range("H3").Value = "2016-03-24 19:27:22"
$de.Range("J3").Value = "Include"

# Row 4: ffffff7695d231 file
$de.Range("A4").Value = "ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("D4").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf"
$de.Range("E4").Value = "2016-03-24 19:26:54"
$de.Range("F4").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.md"
$de.Range("G4").Value = "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf"
$de.Range("H4").Value = "2016-03-24 19:27:22"
$de.Range("J4").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43be2a4a8b012feba1bb9040003eaf8676f25017/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.de-de.xlf", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/43be2a4a8b012feba1bb9040003eaf8676f25017/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.de-de.xlf", "", "", "c9f7cf44-7eb7-4e7a-8ac2-e5219c72462e.84358e14b9ab4479cd9ec144f16eda0220995e73.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md", "", "", "ffff86e48609-e21e-4e4b-a9b8-89bc6ca21893.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee4e668dcec1786f68e46acfd97b103bf196a60f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1736cfa06673407340d5cb88be4fd3185d5d78fa/e2e/ee944fac-2a76-4622-b02c-316995cd2ba5.md", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a5838e46b1520efba87c6fd658537673be0d32be/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf")

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8a63b8f27b283c81a71eeadd01e640d22dda22be/e2e/ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md", "", "", "ffffff7695d231-d2e2-4c8e-ab04-b446f3536747.md")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee4e668dcec1786f68e46acfd97b103bf196a60f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1736cfa06673407340d5cb88be4fd3185d5d78fa/e2e/ee944fac-2a76-4622-b02c-316995cd2ba5.md", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.md")
$de.Hyperlinks.Add($de.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a5838e46b1520efba87c6fd658537673be0d32be/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf", "", "", "ee944fac-2a76-4622-b02c-316995cd2ba5.6b9a001b6dfe4910feb4f3245bde016bc1e64762.de-de.xlf")

Write-Output "Localization status report regenerated for handback."
